$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A15").Value = "VIBE Lounge & coctail bar"
$ws.Range("B15").Value = "Ulice Masarykova 3125/18, 400 01 Ústí nad Labem-centrum"
$ws.Range("C15").Value = "50.66295701796333"
$ws.Range("D15").Value = "14.034403534402813"
$ws.Range("E15").Value = "vibe_usti"

$ws.Range("F1").Value = "Web"
$ws.Range("G1").Value = "Telefon"

$ws.Range("B14").Select()
